# Add a new data row (row 7) to the "Test" worksheet, mirroring the
# pattern established by the existing rows (row 2 -> row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (border/fill/font/cell-style) from row 6 down onto
# the new row 7 first, so the new cells line up with the existing table
# styling (normal border style for A:C, hyperlink style for D).
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Fill in the new row's values, following on from Vandana05/NAir05/.../Test@127.
$ws.Range("A7").Value = "Vandana06"
$ws.Range("B7").Value = "NAir06"
$ws.Range("C7").Value = 7829844600
$ws.Range("D7").Value = "Test@128"

# Give D7 a mailto hyperlink, same as D2:D6 (those all point at
# "mailto:Test@123" with a "Test@123" display label too).
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:Test@123", "", "", "Test@123")

# Adding the hyperlink re-applies Excel's built-in "Hyperlink" cell style,
# which can disturb the cell's formatting - restore it so D7 keeps the same
# look as the rest of column D, then put the real password text back.
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("D7").Value = "Test@128"

# Match the selection left behind by the edit.
$ws.Range("D11").Select()
